$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 91
$ws.Range("K91").Value = "['Belgium', 3, 1, 2]"
$ws.Range("L91").Value = "['Czech Republic', 1, -1, 2]"
$ws.Range("M91").Value = "['Netherlands', 'Belgium', 'Slovenia', 'Hungary']"
$ws.Range("N91").Value = "['Croatia', 'Czech Republic']"
$ws.Range("O91").Value = "[]"

# Row 92
$ws.Range("I92").Value = "['Slovenia', 3, 0, 2]"
$ws.Range("M92").Value = "['Belgium', 'Netherlands', 'Slovenia', 'Hungary']"

# Row 93
$ws.Range("I93").Value = "['Slovenia', 3, 0, 2]"
$ws.Range("M93").Value = "['Netherlands', 'Ukraine', 'Slovenia', 'Hungary']"

# Row 94
$ws.Range("I94").Value = "['Slovenia', 3, 0, 2]"
$ws.Range("M94").Value = "['Netherlands', 'Slovakia', 'Slovenia', 'Hungary']"

# Row 95
$ws.Range("I95").Value = "['Slovenia', 3, 0, 2]"
$ws.Range("M95").Value = "['Netherlands', 'Slovakia', 'Slovenia', 'Hungary']"

# Row 96
$ws.Range("I96").Value = "['Slovenia', 3, 0, 2]"
$ws.Range("M96").Value = "['Netherlands', 'Slovakia', 'Georgia', 'Slovenia']"
$ws.Range("N96").Value = "['Hungary', 'Croatia']"

# Row 97
$ws.Range("I97").Value = "['Slovenia', 3, 0, 2]"
$ws.Range("M97").Value = "['Netherlands', 'Slovakia', 'Georgia', 'Slovenia']"
$ws.Range("N97").Value = "['Hungary', 'Croatia']"

# Row 98
$ws.Range("I98").Value = "['Slovenia', 3, 0, 2]"
$ws.Range("M98").Value = "['Netherlands', 'Georgia', 'Slovakia', 'Slovenia']"
$ws.Range("N98").Value = "['Hungary', 'Croatia']"

# Row 99
$ws.Range("I99").Value = "['Slovenia', 3, 0, 2]"
$ws.Range("M99").Value = "['Netherlands', 'Georgia', 'Slovakia', 'Slovenia']"
$ws.Range("N99").Value = "['Hungary', 'Croatia']"

# Row 100
$ws.Range("I100").Value = "['Slovenia', 3, 0, 2]"
$ws.Range("M100").Value = "['Netherlands', 'Georgia', 'Slovakia', 'Slovenia']"
$ws.Range("N100").Value = "['Hungary', 'Croatia']"
